# Applies the commit's data refresh: Total Assets (52,630,568.25) now flows
# through the various compliance sheets (40Act, IRS, Illiquid, 12d1, 12d2,
# 12d3), plus a couple of column-width tweaks on the 40Act and IRS sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 40Act_Diversification
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("40Act_Diversification")

$ws.Range("H2").Value = 52630568.25
$ws.Range("J2").Value = 42330135.8
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = "None"
$ws.Range("N2").Value = "IBM, JNJ, IBM, JNJ, MRK, MRK"
$ws.Range("O2").Value = 0.2853640423916883
$ws.Range("P2").Value = 0.7146359576083117
$ws.Range("Q2").Value = "(CSCO, 0, 3.59%, 0.00%), (CSCO, 0, 3.57%, 0.00%), (KO, 0, 3.44%, 0.00%), (MCD, 0, 3.43%, 0.00%), (MCD, 0, 3.43%, 0.00%), (CVX, 0, 3.43%, 0.00%), (KO, 0, 3.43%, 0.00%), (CVX, 0, 3.42%, 0.00%), (PG, 0, 3.39%, 0.00%), (PG, 0, 3.39%, 0.00%), (AMGN, 0, 3.27%, 0.00%), (AMGN, 0, 3.27%, 0.00%), (VZ, 0, 3.09%, 0.00%), (VZ, 0, 3.09%, 0.00%), (CSCO, 0, 2.27%, 0.00%), (KO, 0, 1.85%, 0.00%), (MCD, 0, 1.64%, 0.00%), (CVX, 0, 1.46%, 0.00%), (AMGN, 0, 1.31%, 0.00%), (VZ, 0, 0.99%, 0.00%), (PG, 0, 0.71%, 0.00%)"
$ws.Range("T2").Value = 0.1888965223990657

# Column T (20th column) widens from 12 to 21 characters. The COM
# ColumnWidth property is offset from the saved OOXML <col width> by a
# constant 5/6 character in this engine, so back that out to land exactly
# on the target width.
$ws.Columns.Item(20).ColumnWidth = 21 - (5/6)

# ---------------------------------------------------------------------
# IRS_Diversification
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("IRS_Diversification")

$ws.Range("H2").Value = 52630568.25
$ws.Range("I2").Value = 55442243.79000001
$ws.Range("K2").Value = 2631528.4125
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "None"
$ws.Range("O2").Value = "CSCO (4.45%)"
$ws.Range("P2").Value = "CSCO (4.45%)"

$ws.Columns.Item(14).ColumnWidth = 18 - (5/6)

# ---------------------------------------------------------------------
# Illiquid
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Illiquid")
$ws.Range("C2").Value = 52630568.25
$ws.Range("F2").Value = 0.3511409537175195

# ---------------------------------------------------------------------
# 12d1_Other_Investment_Companies
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("12d1_Other_Investment_Companies")
$ws.Range("C2").Value = 52630568.25

# ---------------------------------------------------------------------
# 12d2_Insurance_Companies
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("12d2_Insurance_Companies")
$ws.Range("D2").Value = 52630568.25

# ---------------------------------------------------------------------
# 12d3_Securities_Business
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("12d3_Securities_Business")
$ws.Range("J2").Value = 52630568.25
